# Update generated statistics (view counts) in both the "展览" and
# "全部类型" sheets, matching the source data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 15256
$ws1.Range("F14").Value = 6040
$ws1.Range("F26").Value = 4971
$ws1.Range("F34").Value = 256

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 15256
$ws4.Range("F15").Value = 6040
$ws4.Range("F27").Value = 4971
$ws4.Range("F36").Value = 256
